$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "26.876.28"
Set-TextValue $ws.Range("E2") "  +0.47%  "

Set-TextValue $ws.Range("D3") "1.643.30"
Set-TextValue $ws.Range("E3") "  +0.10%  "

Set-TextValue $ws.Range("E4") "  -0.88%  "

Set-TextValue $ws.Range("D5") "216.84"
Set-TextValue $ws.Range("E5") "  -0.69%  "

Set-TextValue $ws.Range("E6") "  +0.69%  "

Set-TextValue $ws.Range("E7") "  -0.68%  "

Set-TextValue $ws.Range("E8") "  +1.07%  "

Set-TextValue $ws.Range("E9") "  -0.34%  "

Set-TextValue $ws.Range("D10") "19.84"
Set-TextValue $ws.Range("E10") "  +4.05%  "

Set-TextValue $ws.Range("E11") "  +0.00%  "

Set-TextValue $ws.Range("D12") "1.873.49"
Set-TextValue $ws.Range("E12") "  +0.20%  "

Set-TextValue $ws.Range("D13") "1.637.89"
Set-TextValue $ws.Range("E13") "  +0.77%  "

Set-TextValue $ws.Range("D14") "4.13"
Set-TextValue $ws.Range("E14") "  -0.09%  "

Set-TextValue $ws.Range("E15") "  +0.64%  "

Set-TextValue $ws.Range("D16") "66.37"
Set-TextValue $ws.Range("E16") "  +2.70%  "

Set-TextValue $ws.Range("D17") "26.889.62"
Set-TextValue $ws.Range("E17") "  +0.63%  "

Set-TextValue $ws.Range("E18") "  +0.80%  "

Set-TextValue $ws.Range("D19") "218.19"
Set-TextValue $ws.Range("E19") "  +3.25%  "

Set-TextValue $ws.Range("E20") "  -0.79%  "

Set-TextValue $ws.Range("E21") "  +1.18%  "

Set-TextValue $ws.Range("E22") "  +7.19%  "

Set-TextValue $ws.Range("D23") "2.44"
Set-TextValue $ws.Range("E23") "  +5.38%  "

Set-TextValue $ws.Range("E24") "  -0.90%  "

Set-TextValue $ws.Range("D25") "145.62"
Set-TextValue $ws.Range("E25") "  -1.04%  "

Set-TextValue $ws.Range("E26") "  -0.85%  "

Set-TextValue $ws.Range("D27") "7.37"
Set-TextValue $ws.Range("E27") "  +4.19%  "

Set-TextValue $ws.Range("E28") "  +0.44%  "

Set-TextValue $ws.Range("D29") "15.86"
Set-TextValue $ws.Range("E29") "  +1.88%  "

Set-TextValue $ws.Range("D30") "0.0512"
Set-TextValue $ws.Range("E30") "  +2.19%  "

Set-TextValue $ws.Range("D31") "1.19"
Set-TextValue $ws.Range("E31") "  +0.03%  "

Set-TextValue $ws.Range("E32") "  +0.37%  "

Set-TextValue $ws.Range("E33") "  +0.43%  "

Set-TextValue $ws.Range("E34") "  +2.07%  "

Set-TextValue $ws.Range("D35") "2.44"
Set-TextValue $ws.Range("E35") "  -0.69%  "

Set-TextValue $ws.Range("D36") "1.245.46"
Set-TextValue $ws.Range("E36") "  -2.29%  "

Set-TextValue $ws.Range("D37") "0.0174"
Set-TextValue $ws.Range("E37") "  -0.28%  "

Set-TextValue $ws.Range("D38") "0.539"
Set-TextValue $ws.Range("E38") "  +2.22%  "

Set-TextValue $ws.Range("D39") "0.835"
Set-TextValue $ws.Range("E39") "  +3.56%  "

Set-TextValue $ws.Range("E40") "  -0.68%  "

Set-TextValue $ws.Range("D41") "0.811"
Set-TextValue $ws.Range("E41") "  +0.90%  "

Set-TextValue $ws.Range("D42") "5.36"
Set-TextValue $ws.Range("E42") "  +1.80%  "

Set-TextValue $ws.Range("D43") "1.786.67"
Set-TextValue $ws.Range("E43") "  +0.44%  "

Set-TextValue $ws.Range("E44") "  -3.55%  "

Set-TextValue $ws.Range("D45") "60.97"
Set-TextValue $ws.Range("E45") "  +1.20%  "

Set-TextValue $ws.Range("D46") "91.45"
Set-TextValue $ws.Range("E46") "  +0.01%  "

Set-TextValue $ws.Range("D47") "1.59"
Set-TextValue $ws.Range("E47") "  +0.73%  "

Set-TextValue $ws.Range("E48") "  +1.23%  "

Set-TextValue $ws.Range("E49") "  -1.14%  "

Set-TextValue $ws.Range("E50") "  +1.39%  "

Set-TextValue $ws.Range("D51") "7.57"
Set-TextValue $ws.Range("E51") "  +0.50%  "
